# Generate Report for Handoff
# Adds two new file entries (3beaf1c6-... and 7d85ea15-...) to the
# localization-status workbook: one new row per file on the "Overview"
# sheet, and one new row per file on each of the "zh-cn" / "de-de" sheets.

$wb = $excel.ActiveWorkbook

$mdBase    = "https://github.com/OpenLocalizationTest/oltest/blob/a220446e4e152ff1241e014fcc6a2d6817b0de9b/e2e/"
$zhcnBase  = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/b8b768f68cd8a26c8fce7f19d9e11e310f05198f/ol-handoff/OpenLocalizationTestOrg/oltest-zhcn-fly/xinjiang/ht/"
$dedeBase  = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/50e062b086f59a3a861133384f781a45f543dbf5/ol-handoff/OpenLocalizationTestOrg/oltest-dede-fly/xinjiang/ht/"

# The two files being handed off, in the order they appear in the diff.
$file1Uuid = "3beaf1c6-3919-4854-95be-e4a9e92a8aea"
$file1Hash = "834ba30845ff8138eab71e1f25b4dd9af76f9789"
$file2Uuid = "7d85ea15-a875-4e90-8644-c5ff20aff989"
$file2Hash = "2b549f467ab28985f718c39a1f9e990ae59c0218"

$status       = "Ready for handoff"
$handoffDate  = "2016-03-22 11:50:36"
$zhHandoffDt  = "2016-03-22 11:50:28"
$deHandoffDt  = "2016-03-22 11:50:36"
$handbackDt   = "0001-01-01 00:00:00"
$reason       = "Include"
$ext          = ".md"

# ---------------------------------------------------------------------
# Sheet "Overview": File Name | zh-cn | de-de | Latest Handoff Date
# ---------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

$overviewRows = @(
    @{ Row = 4; Uuid = $file1Uuid },
    @{ Row = 5; Uuid = $file2Uuid }
)

foreach ($entry in $overviewRows) {
    $r = $entry.Row
    $uuid = $entry.Uuid
    $mdName = "$uuid.md"

    $aCell = $wsOverview.Cells.Item($r, 1)
    $aCell.Value = $mdName
    $wsOverview.Hyperlinks.Add($aCell, "$mdBase$mdName", $null, $null, $mdName)

    $wsOverview.Cells.Item($r, 2).Value = $status
    $wsOverview.Cells.Item($r, 3).Value = $status
    $wsOverview.Cells.Item($r, 4).Value = $handoffDate
}

# ---------------------------------------------------------------------
# Sheets "zh-cn" / "de-de":
# Source File Name | File Extension | Status | Latest Handoff File |
# Latest Handoff Datetime | Latest Target File | Latest Handback File |
# Latest Handback DateTime | Reference Tokens | Handoff Reason | ...
# ---------------------------------------------------------------------
$langSheets = @(
    @{ Name = "zh-cn"; XlfBase = $zhcnBase; Lang = "zh-cn"; HandoffDt = $zhHandoffDt },
    @{ Name = "de-de"; XlfBase = $dedeBase; Lang = "de-de"; HandoffDt = $deHandoffDt }
)

$fileRows = @(
    @{ Row = 4; Uuid = $file1Uuid; Hash = $file1Hash },
    @{ Row = 5; Uuid = $file2Uuid; Hash = $file2Hash }
)

foreach ($langInfo in $langSheets) {
    $ws = $wb.Worksheets.Item($langInfo.Name)

    foreach ($entry in $fileRows) {
        $r = $entry.Row
        $uuid = $entry.Uuid
        $hash = $entry.Hash
        $mdName = "$uuid.md"
        $xlfName = "$uuid.$hash.$($langInfo.Lang).xlf"

        $aCell = $ws.Cells.Item($r, 1)
        $aCell.Value = $mdName
        $ws.Hyperlinks.Add($aCell, "$mdBase$mdName", $null, $null, $mdName)

        $ws.Cells.Item($r, 2).Value = $ext
        $ws.Cells.Item($r, 3).Value = $status

        $dCell = $ws.Cells.Item($r, 4)
        $dCell.Value = $xlfName
        $ws.Hyperlinks.Add($dCell, "$($langInfo.XlfBase)$xlfName", $null, $null, $xlfName)

        $ws.Cells.Item($r, 5).Value = $langInfo.HandoffDt
        $ws.Cells.Item($r, 8).Value = $handbackDt
        $ws.Cells.Item($r, 10).Value = $reason
    }
}

Write-Output "Report rows for handoff added."
